$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1390.1111
$ws.Range("I33").Value = 1533.8
$ws.Range("K33").Value = 1533.8
$ws.Range("M33").Value = -1304.8

$ws.Range("H52").Value = 1000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 1000
$ws.Range("K52").Value = 0
$ws.Range("M52").Value = 3000
$ws.Range("N52").Value = -3320
$ws.Range("L52").ClearContents()

$ws.Range("H74").Value = 7280.625
$ws.Range("I74").Value = 7699.1
$ws.Range("K74").Value = 7699.1
$ws.Range("M74").Value = -6763.1

$ws.Range("H77").Value = 7280.625
$ws.Range("I77").Value = 7699.1
$ws.Range("K77").Value = 38495.5
$ws.Range("M77").Value = -33815.5

$ws.Range("H112").Value = 1986.5834
$ws.Range("J112").Value = 1928
$ws.Range("L112").Value = 5784
$ws.Range("N112").Value = -8000

$ws.Range("H129").Value = 1199.2307
$ws.Range("I129").Value = 698.875
$ws.Range("J129").Value = 1999.8
$ws.Range("K129").Value = 2096.625
$ws.Range("L129").Value = 5999.4
$ws.Range("M129").Value = 2903.375
$ws.Range("N129").Value = -15999.4

$ws.Range("H132").Value = 22226902
$ws.Range("I132").Value = 32262866
$ws.Range("K132").Value = 96788598
$ws.Range("M132").Value = -96786068

$ws.Range("H141").Value = 6998.3335
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13586.7
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H45").Value = 2772.25
$ws.Range("I45").Value = 1896.7
$ws.Range("K45").Value = 1896.7
$ws.Range("M45").Value = -1519.7

$ws.Range("H46").Value = 33810.332
$ws.Range("I46").Value = 21019.8
$ws.Range("J46").Value = 49798.5
$ws.Range("K46").Value = 21019.8
$ws.Range("L46").Value = 49798.5
$ws.Range("M46").Value = -20700.8
$ws.Range("N46").Value = -50436.5

$ws.Range("H92").Value = 84839.375
$ws.Range("J92").Value = 84839.375
$ws.Range("L92").Value = 84839.375
$ws.Range("N92").Value = -89831.375

$ws.Range("H122").Value = 3687.842
$ws.Range("I122").Value = 3782.3057
$ws.Range("J122").Value = 1987.5
$ws.Range("K122").Value = 11346.9171
$ws.Range("L122").Value = 5962.5
$ws.Range("M122").Value = -8896.917099999999
$ws.Range("N122").Value = -10862.5

$ws.Range("H132").Value = 5987.811
$ws.Range("I132").Value = 2251.6667
$ws.Range("K132").Value = 6755.000100000001
$ws.Range("M132").Value = -4225.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 24499.5
$ws.Range("J103").Value = 24499.5
$ws.Range("L103").Value = 24499.5
$ws.Range("N103").Value = -26843.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 162.35715
$ws.Range("I7").Value = 164.7
$ws.Range("J7").Value = 156.5
$ws.Range("K7").Value = 164.7
$ws.Range("L7").Value = 156.5
$ws.Range("M7").Value = -51.69999999999999
$ws.Range("N7").Value = -382.5

$ws.Range("H16").Value = 1822.8182
$ws.Range("I16").Value = 1822.8182
$ws.Range("K16").Value = 1822.8182
$ws.Range("M16").Value = -1535.8182

$ws.Range("H74").Value = 54750
$ws.Range("J74").Value = 54750
$ws.Range("L74").Value = 54750
$ws.Range("N74").Value = -56498

$ws.Range("H77").Value = 54750
$ws.Range("J77").Value = 54750
$ws.Range("L77").Value = 164250
$ws.Range("N77").Value = -172986

$ws.Range("H92").Value = 78935.375
$ws.Range("J92").Value = 78935.375
$ws.Range("L92").Value = 78935.375
$ws.Range("N92").Value = -83927.375

$ws.Range("H94").Value = 2800
$ws.Range("I94").Value = 300
$ws.Range("J94").Value = 3300
$ws.Range("K94").Value = 300
$ws.Range("L94").Value = 3300
$ws.Range("M94").Value = 151
$ws.Range("N94").Value = -4202

$ws.Range("H96").Value = 10162
$ws.Range("J96").Value = 10162
$ws.Range("L96").Value = 10162
$ws.Range("N96").Value = -15654

$ws.Range("H107").Value = 967.5217
$ws.Range("I107").Value = 503.5
$ws.Range("K107").Value = 503.5
$ws.Range("M107").Value = 1416.5

$ws.Range("H113").Value = 1822.8182
$ws.Range("I113").Value = 1822.8182
$ws.Range("K113").Value = 1822.8182
$ws.Range("M113").Value = 347.1818000000001

$ws.Range("H134").Value = 1528.8334
$ws.Range("I134").Value = 1546.4286
$ws.Range("J134").Value = 1504.2
$ws.Range("K134").Value = 4639.2858
$ws.Range("L134").Value = 4512.6
$ws.Range("M134").Value = -2104.2858
$ws.Range("N134").Value = -9582.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 2759.875
$ws.Range("I18").Value = 1450
$ws.Range("J18").Value = 3196.5
$ws.Range("K18").Value = 4350
$ws.Range("L18").Value = 9589.5
$ws.Range("M18").Value = -4181
$ws.Range("N18").Value = -9927.5

$ws.Range("H28").Value = 5410
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 5410
$ws.Range("K28").Value = 0
$ws.Range("M28").Value = 16230
$ws.Range("N28").Value = -16694
$ws.Range("L28").ClearContents()

$ws.Range("H109").Value = 8428
$ws.Range("I109").Value = 3333
$ws.Range("J109").Value = 8994.111000000001
$ws.Range("K109").Value = 9999
$ws.Range("L109").Value = 26982.333
$ws.Range("M109").Value = -8959
$ws.Range("N109").Value = -29062.333

$ws.Range("H113").Value = 1934.2
$ws.Range("J113").Value = 2614
$ws.Range("L113").Value = 7842
$ws.Range("N113").Value = -12182

$ws.Range("H119").Value = 2316.2
$ws.Range("I119").Value = 2237.4285
$ws.Range("K119").Value = 6712.2855
$ws.Range("M119").Value = -1874.2855

$ws.Range("H131").Value = 1718.931
$ws.Range("J131").Value = 1879.7
$ws.Range("L131").Value = 5639.1
$ws.Range("N131").Value = -15719.1

$ws.Range("H132").Value = 2021.2963
$ws.Range("I132").Value = 981.8333
$ws.Range("J132").Value = 2318.2856
$ws.Range("K132").Value = 8836.4997
$ws.Range("L132").Value = 20864.5704
$ws.Range("M132").Value = -6306.4997
$ws.Range("N132").Value = -25924.5704

$ws.Range("H140").Value = 4574.3887
$ws.Range("I140").Value = 3530.3333
$ws.Range("K140").Value = 10590.9999
$ws.Range("M140").Value = -5410.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H113").Value = 3300.0527
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 4428.364
$ws.Range("I122").Value = 3606
$ws.Range("J122").Value = 5867.5
$ws.Range("K122").Value = 10818
$ws.Range("L122").Value = 17602.5
$ws.Range("M122").Value = -8368
$ws.Range("N122").Value = -22502.5

$ws.Range("H124").Value = 70962.836
$ws.Range("J124").Value = 70962.836
$ws.Range("L124").Value = 70962.836
$ws.Range("N124").Value = -80782.836

$ws.Range("H133").Value = 145617
$ws.Range("J133").Value = 145617
$ws.Range("L133").Value = 145617
$ws.Range("N133").Value = -155737

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5716.25
$ws.Range("I132").Value = 5195.3335
$ws.Range("K132").Value = 15586.0005
$ws.Range("M132").Value = -13056.0005

$ws.Range("H136").Value = 5623.44
$ws.Range("I136").Value = 5708.5
$ws.Range("J136").Value = 4999.6665
$ws.Range("K136").Value = 17125.5
$ws.Range("L136").Value = 14998.9995
$ws.Range("M136").Value = -14575.5
$ws.Range("N136").Value = -20098.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 51392.5
$ws.Range("I81").Value = 80202.92
$ws.Range("K81").Value = 160405.84
$ws.Range("M81").Value = -159344.84

$ws.Range("H84").Value = 51392.5
$ws.Range("I84").Value = 80202.92
$ws.Range("K84").Value = 802029.2
$ws.Range("M84").Value = -796725.2

$ws.Range("H132").Value = 3250
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H136").Value = 3974.7896
$ws.Range("I136").Value = 2790.5625
$ws.Range("K136").Value = 8371.6875
$ws.Range("M136").Value = -5821.6875

